$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace duplicated teacher lists with a single "-" placeholder
$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "-"
$ws.Range("E3").Value = "-"

$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "-"
$ws.Range("E4").Value = "-"

$ws.Range("C6").Value = "-"
$ws.Range("E6").Value = "-"

$ws.Range("C7").Value = "-"
$ws.Range("E7").Value = "-"
